$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Trial 8"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 14.3907253742218
$ws.Range("F9").Value = 13.88457727432251
$ws.Range("G9").Value = 16.58563923835754
$ws.Range("H9").Value = 13.90225052833557

$ws.Range("A10").Value = "Trial 9"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 14.10806488990784
$ws.Range("F10").Value = 13.98515152931213
$ws.Range("G10").Value = 30.99015283584595
$ws.Range("H10").Value = 13.90586185455322

$ws.Range("A11").Value = "Trial 10"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 14.1343719959259
$ws.Range("F11").Value = 14.00995993614197
$ws.Range("G11").Value = 63.32477641105652
$ws.Range("H11").Value = 13.89554309844971

$ws.Range("A12").Value = "Trial 11"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 14.10972428321838
$ws.Range("F12").Value = 13.98259258270264
$ws.Range("G12").Value = 137.5324246883392
$ws.Range("H12").Value = 13.90477609634399

$ws.Range("A13").Value = "Trial 12"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 7.791940450668335
$ws.Range("F13").Value = 7.603684425354004
$ws.Range("G13").Value = 4.725376844406128
$ws.Range("H13").Value = 7.510883331298828

$ws.Range("A14").Value = "Trial 13"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 7.796133756637573
$ws.Range("F14").Value = 7.608644247055054
$ws.Range("G14").Value = 8.80347752571106
$ws.Range("H14").Value = 7.557884216308594

$ws.Range("A15").Value = "Trial 14"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 7.748147487640381
$ws.Range("F15").Value = 7.62830924987793
$ws.Range("G15").Value = 18.12826251983643
$ws.Range("H15").Value = 7.569580793380737

$ws.Range("A16").Value = "Trial 15"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 7.83071756362915
$ws.Range("F16").Value = 7.661117315292358
$ws.Range("G16").Value = 41.8465895652771
$ws.Range("H16").Value = 7.605157136917114

$ws.Range("A17").Value = "Trial 16"
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 10.75392031669617
$ws.Range("F17").Value = 10.64612007141113
$ws.Range("G17").Value = 11.42532658576965
$ws.Range("H17").Value = 10.54792332649231

$ws.Range("A18").Value = "Trial 17"
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 10.74261164665222
$ws.Range("F18").Value = 10.65383958816528
$ws.Range("G18").Value = 28.15839242935181
$ws.Range("H18").Value = 10.49790549278259

$ws.Range("A19").Value = "Trial 18"
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 21.34992527961731
$ws.Range("F19").Value = 21.14037585258484
$ws.Range("G19").Value = 21.50210118293762
$ws.Range("H19").Value = 21.14603114128113

